# Updated BGR model - 2025-08-13 10:18
# Re-rank the "lcoe_class" (column P) values for tied-LCOE resource rows
# on the "solar" and "wind" worksheets.

$wb = $excel.ActiveWorkbook

# ---- solar sheet ----
$wsSolar = $wb.Worksheets.Item("solar")
$wsSolar.Range("P4").Value = 4
$wsSolar.Range("P5").Value = 3
$wsSolar.Range("P6").Value = 2

# ---- wind sheet ----
$wsWind = $wb.Worksheets.Item("wind")
$wsWind.Range("P4").Value = 3
$wsWind.Range("P5").Value = 2

$wsWind.Range("P13").Value = 5
$wsWind.Range("P14").Value = 4

$wsWind.Range("P15").Value = 3
$wsWind.Range("P16").Value = 2
$wsWind.Range("P17").Value = 1

$wsWind.Range("P18").Value = 1
$wsWind.Range("P19").Value = 2
$wsWind.Range("P20").Value = 3

$wsWind.Range("P47").Value = 1
$wsWind.Range("P48").Value = 2
